$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently holds rows 178-218 of "Murcott"/"Clementina" Mandarina
# price data. Three brand-new rows describing a new "Clemenuless" variety
# need to be inserted right before the existing row 178, pushing all of the
# old 178-218 data down to 181-221 (dimension grows from T218 to T221).

$ws.Range("A178:T180").EntireRow.Insert()

# --- New row 178: Clemenuless / Especial ---
$ws.Range("A178").Value() = 4
$ws.Range("B178").Value() = "Feria Lagunitas de Puerto Montt"
$ws.Range("C178").Value() = "Los Lagos"
$ws.Range("D178").Value() = 44722
$ws.Range("E178").Value() = 10
$ws.Range("F178").Value() = "Fruta"
$ws.Range("G178").Value() = 100102
$ws.Range("H178").Value() = "Cítricos"
$ws.Range("I178").Value() = 100102004
$ws.Range("J178").Value() = "Mandarina"
$ws.Range("K178").Value() = "Clemenuless"
$ws.Range("L178").Value() = "Especial"
$ws.Range("M178").Value() = 300
$ws.Range("N178").Value() = 12500
$ws.Range("O178").Value() = 12500
$ws.Range("P178").Value() = 12500
$ws.Range("Q178").Value() = "$/bandeja 10 kilos"
$ws.Range("R178").Value() = "Provincia de Limarí"
$ws.Range("S178").Value() = 1250
$ws.Range("T178").Value() = 10

# --- New row 179: Clemenuless / Primera ---
$ws.Range("A179").Value() = 4
$ws.Range("B179").Value() = "Feria Lagunitas de Puerto Montt"
$ws.Range("C179").Value() = "Los Lagos"
$ws.Range("D179").Value() = 44722
$ws.Range("E179").Value() = 10
$ws.Range("F179").Value() = "Fruta"
$ws.Range("G179").Value() = 100102
$ws.Range("H179").Value() = "Cítricos"
$ws.Range("I179").Value() = 100102004
$ws.Range("J179").Value() = "Mandarina"
$ws.Range("K179").Value() = "Clemenuless"
$ws.Range("L179").Value() = "Primera"
$ws.Range("M179").Value() = 300
$ws.Range("N179").Value() = 10000
$ws.Range("O179").Value() = 10000
$ws.Range("P179").Value() = 10000
$ws.Range("Q179").Value() = "$/bandeja 10 kilos"
$ws.Range("R179").Value() = "Provincia de Limarí"
$ws.Range("S179").Value() = 1000
$ws.Range("T179").Value() = 10

# --- New row 180: Clemenuless / Segunda ---
$ws.Range("A180").Value() = 4
$ws.Range("B180").Value() = "Feria Lagunitas de Puerto Montt"
$ws.Range("C180").Value() = "Los Lagos"
$ws.Range("D180").Value() = 44722
$ws.Range("E180").Value() = 10
$ws.Range("F180").Value() = "Fruta"
$ws.Range("G180").Value() = 100102
$ws.Range("H180").Value() = "Cítricos"
$ws.Range("I180").Value() = 100102004
$ws.Range("J180").Value() = "Mandarina"
$ws.Range("K180").Value() = "Clemenuless"
$ws.Range("L180").Value() = "Segunda"
$ws.Range("M180").Value() = 300
$ws.Range("N180").Value() = 8500
$ws.Range("O180").Value() = 8500
$ws.Range("P180").Value() = 8500
$ws.Range("Q180").Value() = "$/bandeja 10 kilos"
$ws.Range("R180").Value() = "Provincia de Limarí"
$ws.Range("S180").Value() = 850
$ws.Range("T180").Value() = 10
